# Update "Analyse des resultatsV2.xlsx"
#
#  - Fill the H:M numeric-score block (rows 5-9) that sits next to the
#    existing N "Rang" column on the MTH1006, MEC1210 and MEC2115 sheets.
#  - Append a new row 11 (A11:F11) of fractional values on MTH1006 and
#    MEC1210 (MEC2115 keeps its original used range).
#  - Restore the saved view state: per-sheet zoom level, the selected
#    cell on each sheet, and which sheet/tab is active when the file is
#    reopened (MEC2115).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# MTH1006 - H5:M9 scores + new row 11 (A11:F11)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MTH1006")

$mth1006 = @(
    @(5, 2, 5, 5, 5, 5),
    @(5, 5, 5, 5, 2, 5),
    @(2, 5, 2, 5, 1, 5),
    @(4, 5, 4, 5, 1, 5),
    @(5, 1, 5, 5, 4, 2)
)
for ($i = 0; $i -lt 5; $i++) {
    $r = 5 + $i
    for ($j = 0; $j -lt 6; $j++) {
        $ws1.Cells.Item($r, 8 + $j).Value = $mth1006[$i][$j]
    }
}

$row11_1006 = @(0.4, 0.8, 0.4, 1, 0.2, 0.6)
for ($j = 0; $j -lt 6; $j++) {
    $ws1.Cells.Item(11, 1 + $j).Value = $row11_1006[$j]
}

# ---------------------------------------------------------------------------
# MEC1210 - H5:M9 scores + new row 11 (A11:F11)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MEC1210")

$mec1210 = @(
    @(4, 3, 4, 4, 5, 4),
    @(4, 1, 4, 4, 5, 4),
    @(5, 1, 5, 5, 4, 5),
    @(5, 1, 5, 5, 5, 5),
    @(5, 2, 5, 5, 5, 5)
)
for ($i = 0; $i -lt 5; $i++) {
    $r = 5 + $i
    for ($j = 0; $j -lt 6; $j++) {
        $ws2.Cells.Item($r, 8 + $j).Value = $mec1210[$i][$j]
    }
}

$row11_1210 = @(0.2, 0, 0.2, 0, 0, 0.2)
for ($j = 0; $j -lt 6; $j++) {
    $ws2.Cells.Item(11, 1 + $j).Value = $row11_1210[$j]
}

# ---------------------------------------------------------------------------
# MEC2115 - H5:M9 scores only (used range stays A1:N9)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("MEC2115")

$mec2115 = @(
    @(4, 4, 4, 4, 1, 4),
    @(5, 3, 5, 5, 1, 2),
    @(5, 1, 5, 1, 2, 3),
    @(3, 3, 3, 5, 2, 5),
    @(1, 1, 1, 5, 1, 1)
)
for ($i = 0; $i -lt 5; $i++) {
    $r = 5 + $i
    for ($j = 0; $j -lt 6; $j++) {
        $ws3.Cells.Item($r, 8 + $j).Value = $mec2115[$i][$j]
    }
}

# ---------------------------------------------------------------------------
# View state: per-sheet zoom + selection, and the active tab (MEC2115 is
# the sheet visible/selected when the workbook is reopened).
# ---------------------------------------------------------------------------
$ws1.Select()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("F12").Select()

$ws2.Select()
$excel.ActiveWindow.Zoom = 70
$ws2.Range("L13").Select()

$ws3.Select()
$excel.ActiveWindow.Zoom = 70
$ws3.Range("G22").Select()
